$wb = $excel.ActiveWorkbook

# --- Explanation sheet: add "Date" conversion row (row 18) ---
$wsExplanation = $wb.Worksheets.Item("Explanation")
$wsExplanation.Range("A18").Value = "Date"
$wsExplanation.Range("B18").Value = "Date"
$wsExplanation.Range("C18").Value = "Date"
$wsExplanation.Range("D18").Value = "%d/%m/%Y"
$wsExplanation.Range("G18").Value = "Date"
$wsExplanation.Range("H18").Value = "Converts the format in column D into ISO 8601 (%Y-%m-%d), errors will be coded in the output to: ""ERR: original value"""
[void]$wsExplanation.Range("H19").Select()

# --- Conversion_Table sheet: add "Date" conversion row (row 17) ---
$wsConversion = $wb.Worksheets.Item("Conversion_Table")
$wsConversion.Range("A17").Value = "Date"
$wsConversion.Range("B17").Value = "Date"
$wsConversion.Range("C17").Value = "Date"
$wsConversion.Range("D17").Value = "%d/%m/%Y"
[void]$wsConversion.Range("D17").Select()
